$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2: project finish date slips from Jan 22 to Jan 25 (B2 DATEDIF
# recalculates automatically from 23 -> 26 once E2 changes).
# ---------------------------------------------------------------------
$ws.Range("E2").Formula = "=DATE(2024,1,25)"

# ---------------------------------------------------------------------
# Row 5 / Row 6: two new tasks under "Project Management:" — fill in
# the previously-blank rows with real data. Copy number formats first
# so the new cells line up with the rest of the table, then set
# content.
# ---------------------------------------------------------------------
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null

$ws.Range("D4:E4").Copy() | Out-Null
$ws.Range("D5:E5").PasteSpecial(-4122) | Out-Null
$ws.Range("D6:E6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "Project Report"
$ws.Range("C5").Value = "Syed Kaif Ali"
$ws.Range("D5").Formula = "=DATE(2024,1,24)"
$ws.Range("E5").Formula = "=DATE(2024,1,25)"
$ws.Range("B5").Formula = "=E5-D5"

$ws.Range("A6").Value = "Timeline of the Project"
$ws.Range("B6").Value = "Ongoing"
$ws.Range("C6").Value = "Syed Kaif Ali"
$ws.Range("D6").Formula = "=DATE(2023,12,30)"
$ws.Range("E6").Formula = "=DATE(2024,1,1)"

# DATEDIF/subtraction formulas on a date-valued range inherit a date
# number format in this engine (mirrors Excel's own auto-format
# heuristic) - paste the plain "General" look back on top afterwards.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Event Creation block (rows 23-26): contributor re-assigned to the
# newly-added Mustassum Tanvir.
# ---------------------------------------------------------------------
$ws.Range("C23").Value = "Mustassum Tanvir"
$ws.Range("C24").Value = "Mustassum Tanvir"
$ws.Range("C25").Value = "Mustassum Tanvir"
$ws.Range("C26").Value = "Mustassum Tanvir"

# ---------------------------------------------------------------------
# User Front-end block (rows 33-41): fix the duplicated "Open events"
# label (row 34 should be "View upcoming events", shifting the rest of
# the task list down by one) and populate Duration/Contributor/Start/
# Finish for every task - previously only the task name existed.
# ---------------------------------------------------------------------
$ws.Range("C19").Copy() | Out-Null
$ws.Range("C33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B33:C33").Copy() | Out-Null
$ws.Range("B34:C41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("E23").Copy() | Out-Null
$ws.Range("E33:E41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("D23").Copy() | Out-Null
$ws.Range("D34:D41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A34").Value = "View upcoming events"
$ws.Range("A35").Value = "View past events"
$ws.Range("A36").Value = "View favorite events"
$ws.Range("A37").Value = "Remove events from favorites"
$ws.Range("A38").Value = "Buy ticket for an event"
$ws.Range("A39").Value = "Generate QR code"
$ws.Range("A40").Value = "Sort events"
$ws.Range("A41").Value = "View events by categories"
$ws.Range("A42").Clear()

$rows33 = @(33,34,35,36,37,38,39,40,41)
$contributors = @("Ahmad Irfan","Ahmad Irfan","Ahmad Irfan","Abdul Majid","Abdul Majid","Abdul Majid","Ahmad Irfan","Ahmad Irfan","Ahmad Irfan")
$starts = @("2024,1,17","2024,1,17","2024,1,17","2024,1,18","2024,1,18","2024,1,18","2024,1,19","2024,1,19","2024,1,19")
$finishes = @("2024,1,21","2024,1,21","2024,1,21","2024,1,23","2024,1,23","2024,1,23","2024,1,25","2024,1,25","2024,1,25")

for ($i = 0; $i -lt $rows33.Length; $i++) {
    $r = $rows33[$i]
    $ws.Range("C$r").Value = $contributors[$i]
    $ws.Range("D$r").Formula = "=DATE(" + $starts[$i] + ")"
    $ws.Range("E$r").Formula = "=DATE(" + $finishes[$i] + ")"
    $ws.Range("B$r").Formula = "=E$r-D$r"
}

# Same date-format inheritance fix as above, applied to the whole
# B33:B41 column in one shot.
$ws.Range("C19").Copy() | Out-Null
$ws.Range("B33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("B33").Copy() | Out-Null
$ws.Range("B34:B41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Trailing blank row 46 is no longer needed.
# ---------------------------------------------------------------------
$ws.Rows("46").Delete()

# ---------------------------------------------------------------------
# Scroll position / selection left behind by the last editing session.
# ---------------------------------------------------------------------
$ws.Range("G8").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
